{"js": "// Update the two-digit \u00f7 one-digit division prompts in the practice\n// table. Each old expression is unique in the document, so a scoped\n// search-and-replace (matchCase, whole text) is safe and unambiguous.\nconst replacements = [\n  [\"27\u00f78=\", \"66\u00f72=\"],\n  [\"78\u00f75=\", \"97\u00f76=\"],\n  [\"55\u00f75=\", \"57\u00f74=\"],\n  [\"84\u00f72=\", \"91\u00f73=\"],\n  [\"47\u00f75=\", \"31\u00f74=\"],\n  [\"51\u00f78=\", \"72\u00f76=\"],\n  [\"67\u00f79=\", \"87\u00f74=\"],\n  [\"28\u00f77=\", \"30\u00f73=\"],\n  [\"24\u00f72=\", \"21\u00f72=\"],\n  [\"62\u00f72=\", \"28\u00f78=\"],\n  [\"22\u00f76=\", \"57\u00f75=\"],\n  [\"50\u00f77=\", \"24\u00f73=\"],\n  [\"77\u00f77=\", \"95\u00f78=\"],\n  [\"41\u00f76=\", \"92\u00f73=\"],\n  [\"13\u00f79=\", \"76\u00f72=\"],\n  [\"43\u00f76=\", \"93\u00f78=\"],\n  [\"32\u00f74=\", \"59\u00f77=\"],\n  [\"36\u00f79=\", \"61\u00f78=\"],\n  [\"24\u00f74=\", \"34\u00f78=\"],\n  [\"83\u00f78=\", \"82\u00f76=\"],\n  [\"95\u00f77=\", \"14\u00f72=\"],\n  [\"19\u00f78=\", \"38\u00f77=\"],\n  [\"34\u00f79=\", \"71\u00f77=\"],\n  [\"46\u00f74=\", \"83\u00f76=\"],\n  [\"68\u00f74=\", \"16\u00f78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const item of found.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit \u00f7 one-digit division prompts in the practice\n# table. Each old expression is unique in the document, so a simple\n# Find/Replace (whole document, match case) for each pair is safe and\n# unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"27\u00f78=\", \"66\u00f72=\"),\n    @(\"78\u00f75=\", \"97\u00f76=\"),\n    @(\"55\u00f75=\", \"57\u00f74=\"),\n    @(\"84\u00f72=\", \"91\u00f73=\"),\n    @(\"47\u00f75=\", \"31\u00f74=\"),\n    @(\"51\u00f78=\", \"72\u00f76=\"),\n    @(\"67\u00f79=\", \"87\u00f74=\"),\n    @(\"28\u00f77=\", \"30\u00f73=\"),\n    @(\"24\u00f72=\", \"21\u00f72=\"),\n    @(\"62\u00f72=\", \"28\u00f78=\"),\n    @(\"22\u00f76=\", \"57\u00f75=\"),\n    @(\"50\u00f77=\", \"24\u00f73=\"),\n    @(\"77\u00f77=\", \"95\u00f78=\"),\n    @(\"41\u00f76=\", \"92\u00f73=\"),\n    @(\"13\u00f79=\", \"76\u00f72=\"),\n    @(\"43\u00f76=\", \"93\u00f78=\"),\n    @(\"32\u00f74=\", \"59\u00f77=\"),\n    @(\"36\u00f79=\", \"61\u00f78=\"),\n    @(\"24\u00f74=\", \"34\u00f78=\"),\n    @(\"83\u00f78=\", \"82\u00f76=\"),\n    @(\"95\u00f77=\", \"14\u00f72=\"),\n    @(\"19\u00f78=\", \"38\u00f77=\"),\n    @(\"34\u00f79=\", \"71\u00f77=\"),\n    @(\"46\u00f74=\", \"83\u00f76=\"),\n    @(\"68\u00f74=\", \"16\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
